# Error Calculations and Plots
# Apply imputation/de-imputation edits to column D (rows still at their
# original row numbers), then remove the two rows that were dropped from
# the source table (original "RM 232" and "SC 92"), which shifts all the
# following rows up by one/two. Finally apply the remaining imputation
# edits on the rows in their new (post-shift) positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: column D edits on rows 2-24 (row numbers unaffected so far) ---
$ws.Range("D2").Value = -13.5
$ws.Range("D6").Value = ""
$ws.Range("D12").Value = -14.1
$ws.Range("D14").Value = ""
$ws.Range("D20").Value = -14
$ws.Range("D21").Value = -14.3
$ws.Range("D23").Value = ""
$ws.Range("D24").Value = ""

# --- Step 2: remove the "RM 232" row (row 26) and the "SC 92" row
#     (originally row 28, which becomes row 27 once row 26 is gone) ---
$ws.Rows("26:26").Delete()
$ws.Rows("27:27").Delete()

# --- Step 3: column B/D edits on the rows in their new positions ---
$ws.Range("B26").Value = -20.2
$ws.Range("B27").Value = ""
$ws.Range("B28").Value = ""
$ws.Range("B29").Value = -19.5
$ws.Range("B30").Value = -19.7
$ws.Range("B31").Value = ""
$ws.Range("D31").Value = -13.7
$ws.Range("B32").Value = ""
$ws.Range("D33").Value = -14.1
